$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "quarantine"
$ws.Range("C2").Value = "quarantine"
$ws.Range("C3").Value = "quarantine"
$ws.Range("A4").Value = "AIA"
$ws.Range("B4").Value = "Anguilla"
$ws.Range("C4").Value = "can travel"
